$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E held local placeholder filenames (tag1.jpg/tag2.jpg/tag3.jpg); replace
# them with the real image URLs. These URLs are exactly the addresses the
# existing hyperlinks on E2/E3/E4 already point to (rId3/rId1/rId2), so the
# cell text now matches its hyperlink target instead of a stale filename.
$ws.Range("E2").Value = "https://www.kneipp.com/dw/image/v2/BGQM_PRD/on/demandware.static/-/Sites-master-catalog/de_CH/dw5650ffd1/918957_front.png?sw=500&sh=600&sm=fit&sfrm=png"
$ws.Range("E3").Value = "https://juraforum.b-cdn.net/img/lx/33803-280x210_8932"
$ws.Range("E4").Value = "https://d2exd72xrrp1s7.cloudfront.net/www/1e/1eq9jfstgx4d01okpw2wvj32454eh2hv3a-p307462067-full/184e8018e3b?width=2880&crop=false&q=70"

# The three existing hyperlinks (E2, E3, E4) keep pointing at the same targets,
# but no longer carry an explicit display-text override (the cell's own text is
# used for display instead).
foreach ($h in $ws.Hyperlinks) {
    $h.TextToDisplay = ""
}

# Update the saved cursor/selection position.
$null = $ws.Range("D7").Select()
